$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "27.242.43"
$ws.Cells.Item(2,5).Value = "  -3.41%  "

$ws.Cells.Item(3,4).Value = "1.809.49"
$ws.Cells.Item(3,5).Value = "  -3.54%  "

$ws.Cells.Item(4,4).Value = "'1.001"
$ws.Cells.Item(4,5).Value = "  -0.08%  "

$ws.Cells.Item(5,4).Value = "'310.45"
$ws.Cells.Item(5,5).Value = "  -1.71%  "

$ws.Cells.Item(6,4).Value = "'1.001"
$ws.Cells.Item(6,5).Value = "  -0.10%  "

$ws.Cells.Item(7,4).Value = "'0.4211"
$ws.Cells.Item(7,5).Value = "  -2.28%  "

$ws.Cells.Item(8,4).Value = "'0.3559"
$ws.Cells.Item(8,5).Value = "  -3.61%  "

$ws.Cells.Item(9,4).Value = "'0.07117"
$ws.Cells.Item(9,5).Value = "  -4.03%  "

$ws.Cells.Item(10,4).Value = "'0.8509"
$ws.Cells.Item(10,5).Value = "  -3.56%  "

$ws.Cells.Item(11,4).Value = "'20.23"
$ws.Cells.Item(11,5).Value = "  -4.29%  "

$ws.Cells.Item(12,4).Value = "1.782.86"
$ws.Cells.Item(12,5).Value = "  -5.29%  "

$ws.Cells.Item(13,4).Value = "'5.320"
$ws.Cells.Item(13,5).Value = "  -3.01%  "

$ws.Cells.Item(14,4).Value = "'6.386"
$ws.Cells.Item(14,5).Value = "  -3.62%  "

$ws.Cells.Item(15,4).Value = "'0.06885"
$ws.Cells.Item(15,5).Value = "  -1.35%  "

$ws.Cells.Item(16,4).Value = "'1.005"
$ws.Cells.Item(16,5).Value = "  +0.20%  "

$ws.Cells.Item(17,4).Value = "'81.34"
$ws.Cells.Item(17,5).Value = "  +0.26%  "

$ws.Cells.Item(18,4).Value = "'0.000008743"
$ws.Cells.Item(18,5).Value = "  -4.44%  "

$ws.Cells.Item(19,4).Value = "'1.001"
$ws.Cells.Item(19,5).Value = "  -0.08%  "

$ws.Cells.Item(20,4).Value = "'15.10"
$ws.Cells.Item(20,5).Value = "  -3.19%  "

$ws.Cells.Item(21,4).Value = "27.090.74"
$ws.Cells.Item(21,5).Value = "  -3.94%  "

$ws.Cells.Item(22,4).Value = "'5.093"
$ws.Cells.Item(22,5).Value = "  +0.15%  "

$ws.Cells.Item(23,4).Value = "'10.85"
$ws.Cells.Item(23,5).Value = "  -0.72%  "

$ws.Cells.Item(24,4).Value = "2.008.82"
$ws.Cells.Item(24,5).Value = "  -7.23%  "

$ws.Cells.Item(25,4).Value = "'1.965"
$ws.Cells.Item(25,5).Value = "  -0.51%  "

$ws.Cells.Item(26,4).Value = "'153.80"
$ws.Cells.Item(26,5).Value = "  -0.05%  "

$ws.Cells.Item(27,4).Value = "'18.22"
$ws.Cells.Item(27,5).Value = "  -2.84%  "

$ws.Cells.Item(28,4).Value = "'5.039"
$ws.Cells.Item(28,5).Value = "  -7.01%  "

$ws.Cells.Item(29,4).Value = "'113.35"
$ws.Cells.Item(29,5).Value = "  -3.65%  "

$ws.Cells.Item(30,4).Value = "'1.705"
$ws.Cells.Item(30,5).Value = "  -8.99%  "

$ws.Cells.Item(31,4).Value = "'0.08894"
$ws.Cells.Item(31,5).Value = "  -0.80%  "

$ws.Cells.Item(32,4).Value = "'0.7432"
$ws.Cells.Item(32,5).Value = "  -6.07%  "

$ws.Cells.Item(33,4).Value = "'4.452"
$ws.Cells.Item(33,5).Value = "  -5.54%  "

$ws.Cells.Item(34,4).Value = "'2.913"
$ws.Cells.Item(34,5).Value = "  -1.93%  "

$ws.Cells.Item(35,4).Value = "'1.103"
$ws.Cells.Item(35,5).Value = "  -6.35%  "

$ws.Cells.Item(36,5).Value = "  -0.04%  "

$ws.Cells.Item(37,4).Value = "'1.069"
$ws.Cells.Item(37,5).Value = "  -5.41%  "

$ws.Cells.Item(38,4).Value = "'0.05195"
$ws.Cells.Item(38,5).Value = "  -4.82%  "

$ws.Cells.Item(39,4).Value = "'0.01900"
$ws.Cells.Item(39,5).Value = "  -3.18%  "

$ws.Cells.Item(40,4).Value = "'0.1637"
$ws.Cells.Item(40,5).Value = "  -3.28%  "

$ws.Cells.Item(41,4).Value = "'2.718"
$ws.Cells.Item(41,5).Value = "  -6.11%  "

$ws.Cells.Item(42,4).Value = "'0.4962"
$ws.Cells.Item(42,5).Value = "  -4.02%  "

$ws.Cells.Item(43,4).Value = "'6.286"
$ws.Cells.Item(43,5).Value = "  -8.41%  "

$ws.Cells.Item(44,4).Value = "'8.172"
$ws.Cells.Item(44,5).Value = "  -5.32%  "

$ws.Cells.Item(45,4).Value = "'105.20"
$ws.Cells.Item(45,5).Value = "  -0.54%  "

$ws.Cells.Item(46,4).Value = "'10.25"
$ws.Cells.Item(46,5).Value = "  -2.89%  "

$ws.Cells.Item(47,5).Value = "  -0.05%  "

$ws.Cells.Item(48,4).Value = "'0.06392"
$ws.Cells.Item(48,5).Value = "  -2.90%  "

$ws.Cells.Item(49,4).Value = "'0.4557"
$ws.Cells.Item(49,5).Value = "  -4.38%  "

$ws.Cells.Item(50,5).Value = "  -3.64%  "

$ws.Cells.Item(51,4).Value = "'62.87"
$ws.Cells.Item(51,5).Value = "  -4.02%  "
